# GDE-8270: pushed updated datasets
#
# The "TL_FX_04_Files" dataset refresh: the mid-rate quoted for EUR/USD
# moved from 3.55 to 3.21, and the whole sheet's effectiveDate column
# (column H, every data row) rolls forward from 2020-04-25 to 2021-10-28.
# The AUD/USD mid-rate (row 15, currently 1.11) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- EUR/USD midRate (row 8, column F): 3.55 -> 3.21 -----------------------
# Leading apostrophe keeps it text (matches the sheet's existing quotePrefix
# text style) instead of letting Excel coerce it to a number.
$ws.Cells.Item(8, 6).Value = "'3.21"

# --- effectiveDate column (H2:H20): 2020-04-25 -> 2021-10-28 ---------------
# Same text-quoting trick: these cells are quote-prefixed text, not real
# dates, so a plain date-looking string must not be auto-converted.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 20 }
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "'2021-10-28"
}

# --- restore the active selection left by the editor ------------------------
$ws.Activate() | Out-Null
$ws.Range("N18").Select() | Out-Null
